$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.978890784935764, 50.054360842727824]"
$ws.Range("T2").Value = "[49.976475402502885, 50.03008654786464]"
$ws.Range("L3").Value = "[49.97490876966586, 50.087721374057054]"
$ws.Range("T3").Value = "[49.95988625299445, 50.03313003720456]"
